# Applies the "Updated cryptos list" price/volume refresh (and the
# Toncoin/InjectiveProtocol row swap at rows 29-30) described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.390.07"
$ws.Range("E2").Value = "  +0.84%  "

# Row 3
$ws.Range("D3").Value = "2.368.91"
$ws.Range("E3").Value = "  +2.68%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'310.17"
$ws.Range("E5").Value = "  +0.00%  "

# Row 6
$ws.Range("D6").Value = "'104.45"
$ws.Range("E6").Value = "  +3.86%  "

# Row 7
$ws.Range("D7").Value = "'0.524"
$ws.Range("E7").Value = "  -2.46%  "

# Row 8
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("D9").Value = "'0.517"
$ws.Range("E9").Value = "  +0.67%  "

# Row 10
$ws.Range("D10").Value = "'36.16"
$ws.Range("E10").Value = "  +0.27%  "

# Row 11
$ws.Range("D11").Value = "'53.17"
$ws.Range("E11").Value = "  +2.11%  "

# Row 12
$ws.Range("D12").Value = "'0.0812"
$ws.Range("E12").Value = "  -1.00%  "

# Row 13
$ws.Range("E13").Value = "  -0.54%  "

# Row 14
$ws.Range("D14").Value = "'7.00"
$ws.Range("E14").Value = "  -0.69%  "

# Row 15
$ws.Range("D15").Value = "2.738.30"
$ws.Range("E15").Value = "  +2.73%  "

# Row 16
$ws.Range("D16").Value = "'15.61"
$ws.Range("E16").Value = "  +4.65%  "

# Row 17
$ws.Range("D17").Value = "2.374.56"
$ws.Range("E17").Value = "  +2.95%  "

# Row 18
$ws.Range("D18").Value = "'0.814"
$ws.Range("E18").Value = "  +1.28%  "

# Row 19
$ws.Range("D19").Value = "43.352.93"
$ws.Range("E19").Value = "  +0.76%  "

# Row 20
$ws.Range("D20").Value = "'11.99"
$ws.Range("E20").Value = "  -4.54%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0925"
$ws.Range("E21").Value = "  +0.31%  "

# Row 22
$ws.Range("D22").Value = "'6.26"
$ws.Range("E22").Value = "  +3.03%  "

# Row 23
$ws.Range("D23").Value = "'68.38"
$ws.Range("E23").Value = "  +0.27%  "

# Row 24
$ws.Range("D24").Value = "'241.70"
$ws.Range("E24").Value = "  +0.63%  "

# Row 25
$ws.Range("E25").Value = "  +1.95%  "

# Row 26
$ws.Range("D26").Value = "'2.64"
$ws.Range("E26").Value = "  +0.59%  "

# Row 27
$ws.Range("E27").Value = "  -0.06%  "

# Row 28
$ws.Range("D28").Value = "'25.79"
$ws.Range("E28").Value = "  +6.88%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.27"
$ws.Range("E29").Value = "  +7.32%  "

# Row 30
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "'36.98"
$ws.Range("E30").Value = "  -5.27%  "

# Row 31
$ws.Range("D31").Value = "'9.60"
$ws.Range("E31").Value = "  -0.52%  "

# Row 32
$ws.Range("D32").Value = "'162.16"
$ws.Range("E32").Value = "  -3.90%  "

# Row 33
$ws.Range("D33").Value = "'5.28"
$ws.Range("E33").Value = "  -1.13%  "

# Row 34
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.04%  "

# Row 35
$ws.Range("D35").Value = "'18.28"
$ws.Range("E35").Value = "  +3.11%  "

# Row 36
$ws.Range("E36").Value = "  +6.68%  "

# Row 37
$ws.Range("D37").Value = "'3.13"
$ws.Range("E37").Value = "  -0.76%  "

# Row 38
$ws.Range("D38").Value = "'4.72"
$ws.Range("E38").Value = "  +11.35%  "

# Row 39
$ws.Range("E39").Value = "  +0.32%  "

# Row 40
$ws.Range("D40").Value = "'1.94"
$ws.Range("E40").Value = "  +5.38%  "

# Row 41
$ws.Range("E41").Value = "  +0.51%  "

# Row 42
$ws.Range("E42").Value = "  -1.55%  "

# Row 43
$ws.Range("D43").Value = "'2.46"
$ws.Range("E43").Value = "  +7.27%  "

# Row 44
$ws.Range("D44").Value = "'20.35"
$ws.Range("E44").Value = "  +5.50%  "

# Row 45
$ws.Range("D45").Value = "2.000.77"
$ws.Range("E45").Value = "  +1.28%  "

# Row 46
$ws.Range("D46").Value = "'0.0290"
$ws.Range("E46").Value = "  +0.21%  "

# Row 47
$ws.Range("D47").Value = "'3.15"
$ws.Range("E47").Value = "  +4.73%  "

# Row 48
$ws.Range("D48").Value = "'10.41"
$ws.Range("E48").Value = "  +6.28%  "

# Row 49
$ws.Range("D49").Value = "'58.20"
$ws.Range("E49").Value = "  +5.37%  "

# Row 50
$ws.Range("D50").Value = "'2.97"
$ws.Range("E50").Value = "  +0.02%  "

# Row 51
$ws.Range("D51").Value = "2.603.80"
$ws.Range("E51").Value = "  +2.77%  "
